$wb = $excel.ActiveWorkbook

# Duplicate the "PL_retirement_adjustment" sheet and place the copy right
# after "Info" (i.e. immediately before "PL_retirement_adjustment"), then
# rename it and make it the active sheet - mirroring the new
# "PL_students_adjustment" worksheet introduced by the edit.
$source = $wb.Worksheets.Item("PL_retirement_adjustment")
$source.Copy($wb.Worksheets.Item("PL_retirement_adjustment"))

$newSheet = $wb.Worksheets.Item("PL_retirement_adjustment (2)")
$newSheet.Name = "PL_students_adjustment"
$newSheet.Activate()
